$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM recomputation values (ligand-expressing cell counts,
# detection rates, average/total expression values, derived specificities, and the
# downstream edge-weight / specificity columns that depend on them).

$row = 2
$ws.Cells.Item($row, 5).Value2 = 3  # E2
$ws.Cells.Item($row, 6).Value2 = 1  # F2
$ws.Cells.Item($row, 7).Value2 = 1.729584333333333  # G2
$ws.Cells.Item($row, 8).Value2 = 5.188753  # H2
$ws.Cells.Item($row, 9).Value2 = 0.2476387648475193  # I2
$ws.Cells.Item($row, 10).Value2 = 0.2476387648475193  # J2
$ws.Cells.Item($row, 13).Value2 = 14.25737566666667  # M2
$ws.Cells.Item($row, 14).Value2 = 42.772127  # N2
$ws.Cells.Item($row, 15).Value2 = 0.2087950866344732  # O2
$ws.Cells.Item($row, 16).Value2 = 0.2087950866344732  # P2
$ws.Cells.Item($row, 17).Value2 = 24.65933358751456  # Q2
$ws.Cells.Item($row, 18).Value2 = 221.934002287631  # R2
$ws.Cells.Item($row, 19).Value2 = 0.05170575736039172  # S2
$ws.Cells.Item($row, 20).Value2 = 0.05170575736039173  # T2

$row = 3
$ws.Cells.Item($row, 5).Value2 = 3  # E3
$ws.Cells.Item($row, 6).Value2 = 1  # F3
$ws.Cells.Item($row, 7).Value2 = 1.729584333333333  # G3
$ws.Cells.Item($row, 8).Value2 = 5.188753  # H3
$ws.Cells.Item($row, 9).Value2 = 0.2476387648475193  # I3
$ws.Cells.Item($row, 10).Value2 = 0.2476387648475193  # J3
$ws.Cells.Item($row, 14).Value2 = 87.128332  # N3
$ws.Cells.Item($row, 15).Value2 = 0.4253229592313036  # O3
$ws.Cells.Item($row, 16).Value2 = 0.4253229592313036  # P3
$ws.Cells.Item($row, 17).Value2 = 50.23193267222178  # Q3
$ws.Cells.Item($row, 18).Value2 = 452.087394049996  # R3
$ws.Cells.Item($row, 19).Value2 = 0.1053264522853318  # S3
$ws.Cells.Item($row, 20).Value2 = 0.1053264522853318  # T3

$row = 4
$ws.Cells.Item($row, 5).Value2 = 3  # E4
$ws.Cells.Item($row, 6).Value2 = 1  # F4
$ws.Cells.Item($row, 7).Value2 = 1.729584333333333  # G4
$ws.Cells.Item($row, 8).Value2 = 5.188753  # H4
$ws.Cells.Item($row, 9).Value2 = 0.2476387648475193  # I4
$ws.Cells.Item($row, 10).Value2 = 0.2476387648475193  # J4
$ws.Cells.Item($row, 13).Value2 = 20.11084633333333  # M4
$ws.Cells.Item($row, 14).Value2 = 60.332539  # N4
$ws.Cells.Item($row, 15).Value2 = 0.2945174484164121  # O4
$ws.Cells.Item($row, 16).Value2 = 0.2945174484164122  # P4
$ws.Cells.Item($row, 17).Value2 = 34.78340474820745  # Q4
$ws.Cells.Item($row, 18).Value2 = 313.050642733867  # R4
$ws.Cells.Item($row, 19).Value2 = 0.07293393715188327  # S4
$ws.Cells.Item($row, 20).Value2 = 0.07293393715188329  # T4

$row = 5
$ws.Cells.Item($row, 5).Value2 = 3  # E5
$ws.Cells.Item($row, 6).Value2 = 1  # F5
$ws.Cells.Item($row, 7).Value2 = 1.729584333333333  # G5
$ws.Cells.Item($row, 8).Value2 = 5.188753  # H5
$ws.Cells.Item($row, 9).Value2 = 0.2476387648475193  # I5
$ws.Cells.Item($row, 10).Value2 = 0.2476387648475193  # J5
$ws.Cells.Item($row, 13).Value2 = 4.873057999999999  # M5
$ws.Cells.Item($row, 14).Value2 = 14.619174  # N5
$ws.Cells.Item($row, 15).Value2 = 0.07136450571781097  # O5
$ws.Cells.Item($row, 16).Value2 = 0.07136450571781099  # P5
$ws.Cells.Item($row, 17).Value2 = 8.428364772224667  # Q5
$ws.Cells.Item($row, 18).Value2 = 75.855282950022  # R5
$ws.Cells.Item($row, 19).Value2 = 0.01767261804991244  # S5
$ws.Cells.Item($row, 20).Value2 = 0.01767261804991244  # T5

$row = 6
$ws.Cells.Item($row, 9).Value2 = 0.2307941364328804  # I6
$ws.Cells.Item($row, 10).Value2 = 0.2307941364328804  # J6
$ws.Cells.Item($row, 13).Value2 = 14.25737566666667  # M6
$ws.Cells.Item($row, 14).Value2 = 42.772127  # N6
$ws.Cells.Item($row, 15).Value2 = 0.2087950866344732  # O6
$ws.Cells.Item($row, 16).Value2 = 0.2087950866344732  # P6
$ws.Cells.Item($row, 17).Value2 = 22.98198185508256  # Q6
$ws.Cells.Item($row, 18).Value2 = 206.837836695743  # R6
$ws.Cells.Item($row, 19).Value2 = 0.04818868171123168  # S6
$ws.Cells.Item($row, 20).Value2 = 0.04818868171123169  # T6

$row = 7
$ws.Cells.Item($row, 9).Value2 = 0.2307941364328804  # I7
$ws.Cells.Item($row, 10).Value2 = 0.2307941364328804  # J7
$ws.Cells.Item($row, 14).Value2 = 87.128332  # N7
$ws.Cells.Item($row, 15).Value2 = 0.4253229592313036  # O7
$ws.Cells.Item($row, 16).Value2 = 0.4253229592313036  # P7
$ws.Cells.Item($row, 17).Value2 = 46.81510800450978  # Q7
$ws.Cells.Item($row, 18).Value2 = 421.335972040588  # R7
$ws.Cells.Item($row, 19).Value2 = 0.0981620450808659  # S7
$ws.Cells.Item($row, 20).Value2 = 0.09816204508086591  # T7

$row = 8
$ws.Cells.Item($row, 9).Value2 = 0.2307941364328804  # I8
$ws.Cells.Item($row, 10).Value2 = 0.2307941364328804  # J8
$ws.Cells.Item($row, 13).Value2 = 20.11084633333333  # M8
$ws.Cells.Item($row, 14).Value2 = 60.332539  # N8
$ws.Cells.Item($row, 15).Value2 = 0.2945174484164121  # O8
$ws.Cells.Item($row, 16).Value2 = 0.2945174484164122  # P8
$ws.Cells.Item($row, 17).Value2 = 32.41740389878344  # Q8
$ws.Cells.Item($row, 18).Value2 = 291.756635089051  # R8
$ws.Cells.Item($row, 19).Value2 = 0.06797290017168123  # S8
$ws.Cells.Item($row, 20).Value2 = 0.06797290017168124  # T8

$row = 9
$ws.Cells.Item($row, 9).Value2 = 0.2307941364328804  # I9
$ws.Cells.Item($row, 10).Value2 = 0.2307941364328804  # J9
$ws.Cells.Item($row, 13).Value2 = 4.873057999999999  # M9
$ws.Cells.Item($row, 14).Value2 = 14.619174  # N9
$ws.Cells.Item($row, 15).Value2 = 0.07136450571781097  # O9
$ws.Cells.Item($row, 16).Value2 = 0.07136450571781099  # P9
$ws.Cells.Item($row, 17).Value2 = 7.855059244640666  # Q9
$ws.Cells.Item($row, 18).Value2 = 70.695533201766  # R9
$ws.Cells.Item($row, 19).Value2 = 0.01647050946910154  # S9
$ws.Cells.Item($row, 20).Value2 = 0.01647050946910154  # T9

$row = 10
$ws.Cells.Item($row, 7).Value2 = 2.743651333333334  # G10
$ws.Cells.Item($row, 8).Value2 = 8.230954000000001  # H10
$ws.Cells.Item($row, 9).Value2 = 0.3928310486309039  # I10
$ws.Cells.Item($row, 10).Value2 = 0.3928310486309038  # J10
$ws.Cells.Item($row, 13).Value2 = 14.25737566666667  # M10
$ws.Cells.Item($row, 14).Value2 = 42.772127  # N10
$ws.Cells.Item($row, 15).Value2 = 0.2087950866344732  # O10
$ws.Cells.Item($row, 16).Value2 = 0.2087950866344732  # P10
$ws.Cells.Item($row, 17).Value2 = 39.11726775768423  # Q10
$ws.Cells.Item($row, 18).Value2 = 352.0554098191581  # R10
$ws.Cells.Item($row, 19).Value2 = 0.08202119283160053  # S10
$ws.Cells.Item($row, 20).Value2 = 0.08202119283160053  # T10

$row = 11
$ws.Cells.Item($row, 7).Value2 = 2.743651333333334  # G11
$ws.Cells.Item($row, 8).Value2 = 8.230954000000001  # H11
$ws.Cells.Item($row, 9).Value2 = 0.3928310486309039  # I11
$ws.Cells.Item($row, 10).Value2 = 0.3928310486309038  # J11
$ws.Cells.Item($row, 14).Value2 = 87.128332  # N11
$ws.Cells.Item($row, 15).Value2 = 0.4253229592313036  # O11
$ws.Cells.Item($row, 16).Value2 = 0.4253229592313036  # P11
$ws.Cells.Item($row, 17).Value2 = 79.68325475430312  # Q11
$ws.Cells.Item($row, 18).Value2 = 717.1492927887281  # R11
$ws.Cells.Item($row, 19).Value2 = 0.1670800640816322  # S11
$ws.Cells.Item($row, 20).Value2 = 0.1670800640816322  # T11

$row = 12
$ws.Cells.Item($row, 7).Value2 = 2.743651333333334  # G12
$ws.Cells.Item($row, 8).Value2 = 8.230954000000001  # H12
$ws.Cells.Item($row, 9).Value2 = 0.3928310486309039  # I12
$ws.Cells.Item($row, 10).Value2 = 0.3928310486309038  # J12
$ws.Cells.Item($row, 13).Value2 = 20.11084633333333  # M12
$ws.Cells.Item($row, 14).Value2 = 60.332539  # N12
$ws.Cells.Item($row, 15).Value2 = 0.2945174484164121  # O12
$ws.Cells.Item($row, 16).Value2 = 0.2945174484164122  # P12
$ws.Cells.Item($row, 17).Value2 = 55.17715035691178  # Q12
$ws.Cells.Item($row, 18).Value2 = 496.594353212206  # R12
$ws.Cells.Item($row, 19).Value2 = 0.1156955981015173  # S12
$ws.Cells.Item($row, 20).Value2 = 0.1156955981015173  # T12

$row = 13
$ws.Cells.Item($row, 7).Value2 = 2.743651333333334  # G13
$ws.Cells.Item($row, 8).Value2 = 8.230954000000001  # H13
$ws.Cells.Item($row, 9).Value2 = 0.3928310486309039  # I13
$ws.Cells.Item($row, 10).Value2 = 0.3928310486309038  # J13
$ws.Cells.Item($row, 13).Value2 = 4.873057999999999  # M13
$ws.Cells.Item($row, 14).Value2 = 14.619174  # N13
$ws.Cells.Item($row, 15).Value2 = 0.07136450571781097  # O13
$ws.Cells.Item($row, 16).Value2 = 0.07136450571781099  # P13
$ws.Cells.Item($row, 17).Value2 = 13.36997207911067  # Q13
$ws.Cells.Item($row, 18).Value2 = 120.329748711996  # R13
$ws.Cells.Item($row, 19).Value2 = 0.02803419361615382  # S13
$ws.Cells.Item($row, 20).Value2 = 0.02803419361615382  # T13

$row = 14
$ws.Cells.Item($row, 5).Value2 = 3  # E14
$ws.Cells.Item($row, 6).Value2 = 1  # F14
$ws.Cells.Item($row, 7).Value2 = 0.8991316666666668  # G14
$ws.Cells.Item($row, 8).Value2 = 2.697395  # H14
$ws.Cells.Item($row, 9).Value2 = 0.1287360500886965  # I14
$ws.Cells.Item($row, 10).Value2 = 0.1287360500886965  # J14
$ws.Cells.Item($row, 13).Value2 = 14.25737566666667  # M14
$ws.Cells.Item($row, 14).Value2 = 42.772127  # N14
$ws.Cells.Item($row, 15).Value2 = 0.2087950866344732  # O14
$ws.Cells.Item($row, 16).Value2 = 0.2087950866344732  # P14
$ws.Cells.Item($row, 17).Value2 = 12.81925794546278  # Q14
$ws.Cells.Item($row, 18).Value2 = 115.373321509165  # R14
$ws.Cells.Item($row, 19).Value2 = 0.02687945473124927  # S14
$ws.Cells.Item($row, 20).Value2 = 0.02687945473124927  # T14

$row = 15
$ws.Cells.Item($row, 5).Value2 = 3  # E15
$ws.Cells.Item($row, 6).Value2 = 1  # F15
$ws.Cells.Item($row, 7).Value2 = 0.8991316666666668  # G15
$ws.Cells.Item($row, 8).Value2 = 2.697395  # H15
$ws.Cells.Item($row, 9).Value2 = 0.1287360500886965  # I15
$ws.Cells.Item($row, 10).Value2 = 0.1287360500886965  # J15
$ws.Cells.Item($row, 14).Value2 = 87.128332  # N15
$ws.Cells.Item($row, 15).Value2 = 0.4253229592313036  # O15
$ws.Cells.Item($row, 16).Value2 = 0.4253229592313036  # P15
$ws.Cells.Item($row, 17).Value2 = 26.11328078834889  # Q15
$ws.Cells.Item($row, 18).Value2 = 235.01952709514  # R15
$ws.Cells.Item($row, 19).Value2 = 0.05475439778347373  # S15
$ws.Cells.Item($row, 20).Value2 = 0.05475439778347373  # T15

$row = 16
$ws.Cells.Item($row, 5).Value2 = 3  # E16
$ws.Cells.Item($row, 6).Value2 = 1  # F16
$ws.Cells.Item($row, 7).Value2 = 0.8991316666666668  # G16
$ws.Cells.Item($row, 8).Value2 = 2.697395  # H16
$ws.Cells.Item($row, 9).Value2 = 0.1287360500886965  # I16
$ws.Cells.Item($row, 10).Value2 = 0.1287360500886965  # J16
$ws.Cells.Item($row, 13).Value2 = 20.11084633333333  # M16
$ws.Cells.Item($row, 14).Value2 = 60.332539  # N16
$ws.Cells.Item($row, 15).Value2 = 0.2945174484164121  # O16
$ws.Cells.Item($row, 16).Value2 = 0.2945174484164122  # P16
$ws.Cells.Item($row, 17).Value2 = 18.08229878176722  # Q16
$ws.Cells.Item($row, 18).Value2 = 162.740689035905  # R16
$ws.Cells.Item($row, 19).Value2 = 0.03791501299133033  # S16
$ws.Cells.Item($row, 20).Value2 = 0.03791501299133034  # T16

$row = 17
$ws.Cells.Item($row, 5).Value2 = 3  # E17
$ws.Cells.Item($row, 6).Value2 = 1  # F17
$ws.Cells.Item($row, 7).Value2 = 0.8991316666666668  # G17
$ws.Cells.Item($row, 8).Value2 = 2.697395  # H17
$ws.Cells.Item($row, 9).Value2 = 0.1287360500886965  # I17
$ws.Cells.Item($row, 10).Value2 = 0.1287360500886965  # J17
$ws.Cells.Item($row, 13).Value2 = 4.873057999999999  # M17
$ws.Cells.Item($row, 14).Value2 = 14.619174  # N17
$ws.Cells.Item($row, 15).Value2 = 0.07136450571781097  # O17
$ws.Cells.Item($row, 16).Value2 = 0.07136450571781099  # P17
$ws.Cells.Item($row, 17).Value2 = 4.381520761303333  # Q17
$ws.Cells.Item($row, 18).Value2 = 39.43368685173  # R17
$ws.Cells.Item($row, 19).Value2 = 0.009187184582643181  # S17
$ws.Cells.Item($row, 20).Value2 = 0.009187184582643183  # T17
